$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.547.90"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.08%  "
$ws.Cells.Item(3, 4).Value = "'1.877.94"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.84%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).Value = "'0.7240"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.28%  "
$ws.Cells.Item(6, 4).Value = "'239.98"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.97%  "
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
$ws.Cells.Item(8, 4).Value = "'0.07852"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -4.00%  "
$ws.Cells.Item(9, 4).Value = "'0.3090"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.77%  "
$ws.Cells.Item(10, 4).Value = "'25.34"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +8.65%  "
$ws.Cells.Item(11, 4).Value = "'0.08238"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.81%  "
$ws.Cells.Item(12, 4).Value = "'1.895.61"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +4.33%  "
$ws.Cells.Item(13, 4).Value = "'0.7273"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.81%  "
$ws.Cells.Item(14, 4).Value = "'5.257"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.62%  "
$ws.Cells.Item(15, 4).Value = "'90.27"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.18%  "
$ws.Cells.Item(16, 4).Value = "'29.576.51"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.17%  "
$ws.Cells.Item(17, 4).Value = "'5.855"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.20%  "
$ws.Cells.Item(18, 4).Value = "'243.42"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.98%  "
$ws.Cells.Item(19, 4).Value = "'0.000007875"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.15%  "
$ws.Cells.Item(20, 5).Value = "  -0.13%  "
$ws.Cells.Item(21, 4).Value = "'2.133.50"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.30%  "
$ws.Cells.Item(22, 4).Value = "'1.000"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.10%  "
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).Value = "'7.794"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +5.40%  "
$ws.Cells.Item(25, 4).Value = "'0.1593"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +9.94%  "
$ws.Cells.Item(26, 4).Value = "'162.80"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.37%  "
$ws.Cells.Item(27, 4).Value = "'9.003"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.53%  "
$ws.Cells.Item(28, 4).Value = "'18.38"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.60%  "
$ws.Cells.Item(29, 4).Value = "'1.949"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.19%  "
$ws.Cells.Item(30, 4).Value = "'1.356"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.80%  "
$ws.Cells.Item(31, 4).Value = "'1.485"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.02%  "
$ws.Cells.Item(32, 4).Value = "'4.357"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.68%  "
$ws.Cells.Item(33, 4).Value = "'4.093"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.18%  "
$ws.Cells.Item(34, 5).Value = "  +0.85%  "
$ws.Cells.Item(35, 4).Value = "'1.201"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.60%  "
$ws.Cells.Item(36, 4).Value = "'0.7202"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.81%  "
$ws.Cells.Item(37, 4).Value = "'0.9990"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.06%  "
$ws.Cells.Item(38, 4).Value = "'2.671"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.02%  "
$ws.Cells.Item(39, 4).Value = "'0.01869"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.13%  "
$ws.Cells.Item(40, 5).Value = "  -0.39%  "
$ws.Cells.Item(41, 4).Value = "'1.185.33"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.71%  "
$ws.Cells.Item(42, 4).Value = "'0.9123"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.83%  "
$ws.Cells.Item(43, 4).Value = "'6.012"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.42%  "
$ws.Cells.Item(44, 4).Value = "'0.4333"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.36%  "
$ws.Cells.Item(45, 4).Value = "'71.95"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.70%  "
$ws.Cells.Item(46, 5).Value = "  +0.16%  "
$ws.Cells.Item(47, 4).Value = "'103.15"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.37%  "
$ws.Cells.Item(48, 4).Value = "'0.5353"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.25%  "
$ws.Cells.Item(49, 4).Value = "'1.783"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.53%  "
$ws.Cells.Item(50, 4).Value = "'2.888"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +5.07%  "
$ws.Cells.Item(51, 4).Value = "'9.249"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.37%  "
